# RND_Todos.xlsx - 2025-06-30 OneDrive sync edit
# Clears the sample row of data from Todos/Updates, removes the single
# Meetings entry, and adds UPDATED_BY / CREATED_BY / CREATED_AT tracking
# columns to the Updates and Meetings sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Todos sheet: wipe the sample row but keep the header + formatting
# ---------------------------------------------------------------
$wsTodos = $wb.Worksheets.Item("Todos")
$wsTodos.Range("A2:H2").ClearContents()
$wsTodos.Rows.Item(2).EntireRow.AutoFit()
$wsTodos.Range("H2").Select()

# ---------------------------------------------------------------
# Updates sheet: add the UPDATED_BY column and wipe the sample row
# ---------------------------------------------------------------
$wsUpdates = $wb.Worksheets.Item("Updates")
$wsUpdates.Range("F1").Value = "UPDATED_BY"
$wsUpdates.Columns.Item(5).ColumnWidth = 13.666666666666666
$wsUpdates.Columns.Item(6).ColumnWidth = 11.333333333333334
$wsUpdates.Range("A2:E2").ClearContents()
$wsUpdates.Range("F1").Select()

# ---------------------------------------------------------------
# Meetings sheet: add CREATED_BY / CREATED_AT columns and remove the
# single logged meeting row
# ---------------------------------------------------------------
$wsMeetings = $wb.Worksheets.Item("Meetings")
$wsMeetings.Range("E1").Value = "CREATED_BY"
$wsMeetings.Range("F1").Value = "CREATED_AT"
$wsMeetings.Columns.Item(5).ColumnWidth = 11.666666666666666
$wsMeetings.Columns.Item(6).ColumnWidth = 11.333333333333334
$wsMeetings.Rows.Item(5).ClearContents()
$wsMeetings.Range("A5").Select()

# Meetings becomes the active sheet/tab after the edit
$wsMeetings.Activate()
